# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 8, pushing the previous
# rows 8 and 9 down to rows 9 and 10 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (existing row 8 -> 9, existing row 9 -> 10)
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with this week's data
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44449
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100112026
$ws.Range("G8").Value = "Haba"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 1300
$ws.Range("K8").Value = 900
$ws.Range("L8").Value = 950
$ws.Range("M8").Value = 925
$ws.Range("N8").Value = "`$/kilo"
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 925
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"
